$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "start"
